$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: ATR case 1 - sample size increased from 73 to 73/73 (all matched)
$ws.Range("D6").Formula = "=73/73"

# Row 16 (E16) text updated first so its new shared string is appended before E9's,
# matching the order in which new strings appear in the target sharedStrings table.
$ws.Range("E16").Value = "random sample of 135 intersections (assumed repsonse distribution of 85%, 5% MOE, 95% LOC)"

# Row 9: ATR case 12 - D9 becomes a computed formula instead of a static value, E9 text updated
$ws.Range("D9").Formula = "=115/121"
$ws.Range("E9").Value = "random sample of 130 (assumed response distribution of 85%, 5% MOE, 95% LOC)"

# Row 16: D16 value changed from 0.85 to 0.9
$ws.Range("D16").Value = 0.9

# Update the selected cell in the sheet view
[void]$ws.Range("D10").Select()
